# Append the 2026-01-09 Kaspa buy to Sheet1 (row 22), matching the existing
# table format: column A is stored as literal text (e.g. "01/02/2026" in
# row 21), not as a date serial. Force text entry (NumberFormat "@") so
# Excel's auto date-recognition doesn't convert the string, then clear the
# formatting again so the new row keeps the sheet's default (unstyled)
# look, same as every other data row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A22").NumberFormat = "@"
$ws.Range("A22").Value = "01/09/2026"
$ws.Range("A22").ClearFormats()

$ws.Range("B22").Value = 1039.873
$ws.Range("C22").Value = 0.04760196677863548
$ws.Range("D22").Value = 50
